# Fix Training Data Issue (#48)
# The "Date" column (BF) held the filename-derived string "6-9-2007-08"
# instead of an actual ISO-style date for the game. Correct it to
# "2008-06-09" for every data row (rows 2-31), keeping the value as text
# (not letting Excel reinterpret "2008-06-09" as a date serial).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "6-9-2007-08"
$newValue = "2008-06-09"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Range("BF$row")
    if ($cell.Value2 -eq $oldValue) {
        # Force text formatting first so Excel doesn't auto-convert the
        # new "YYYY-MM-DD" looking string into a date serial number.
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
    }
}
